# Updated test data for normal load, cable capacitance etc.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices Loop A")

# Row 7: now holds the device data that used to live in row 8 (LI800 / Other),
# but with a brand-new label "LI800 - 1" and refreshed voltage-drop figures.
$ws.Range("A7").Value = "LI800"
$ws.Range("B7").Value = "Other"
$ws.Range("C7").Value = "LI800 - 1"
$ws.Range("E7").Value = 289
$ws.Range("F7").Value = 0.27
$ws.Range("G7").Value = 0.48

# Row 8: now holds the device data that used to live in row 7 (801 CH / Detectors).
$ws.Range("A8").Value = "801 CH"
$ws.Range("B8").Value = "Detectors"
$ws.Range("C8").Value = "801 CH - 3"

# Reflect the new active selection left behind by the edit.
$ws.Range("C7").Select()
